$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_6_2_0"
$ws.Range("B2").Value = 0.2701742343213689
$ws.Range("C2").Value = 0.1131629723332158
$ws.Range("D2").Value = -1.913540015938044
$ws.Range("E2").Value = 0.2121305278927856
$ws.Range("F2").Value = 0.8077020645141602
$ws.Range("G2").Value = 1.56532871723175
$ws.Range("H2").Value = 0.4264284670352936
$ws.Range("I2").Value = 1.029375791549683
$ws.Range("A3").Value = "model_6_2_1"
$ws.Range("B3").Value = 0.2727388848632211
$ws.Range("C3").Value = 0.1134487247300462
$ws.Range("D3").Value = -1.858185881517677
$ws.Range("E3").Value = 0.2152531461068312
$ws.Range("F3").Value = 0.8048638105392456
$ws.Range("G3").Value = 1.564824461936951
$ws.Range("H3").Value = 0.4183267951011658
$ws.Range("I3").Value = 1.025295972824097
$ws.Range("A4").Value = "model_6_2_3"
$ws.Range("B4").Value = 0.2840762013582634
$ws.Range("C4").Value = 0.109720147716511
$ws.Range("D4").Value = -1.506388007714128
$ws.Range("E4").Value = 0.2311319399502455
$ws.Range("F4").Value = 0.792316734790802
$ws.Range("G4").Value = 1.571405649185181
$ws.Range("H4").Value = 0.3668373227119446
$ws.Range("I4").Value = 1.004549860954285
$ws.Range("A5").Value = "model_6_2_2"
$ws.Range("B5").Value = 0.28702344667971
$ws.Range("C5").Value = 0.1407824315618733
$ws.Range("D5").Value = -1.617488467044033
$ws.Range("E5").Value = 0.2474910410603124
$ws.Range("F5").Value = 0.7890548706054688
$ws.Range("G5").Value = 1.516578316688538
$ws.Range("H5").Value = 0.3830980658531189
$ws.Range("I5").Value = 0.9831761717796326
$ws.Range("A6").Value = "model_6_2_12"
$ws.Range("B6").Value = 0.2918265670586877
$ws.Range("C6").Value = 0.1228561128066425
$ws.Range("D6").Value = -1.414509043880367
$ws.Range("E6").Value = 0.2453703497115158
$ws.Range("F6").Value = 0.7837392687797546
$ws.Range("G6").Value = 1.548219561576843
$ws.Range("H6").Value = 0.3533897995948792
$ws.Range("I6").Value = 0.9859468936920166
$ws.Range("A7").Value = "model_6_2_20"
$ws.Range("B7").Value = 0.2925916187066399
$ws.Range("C7").Value = 0.1251258007594539
$ws.Range("D7").Value = -1.419992868707918
$ws.Range("E7").Value = 0.2467048172305097
$ws.Range("F7").Value = 0.782892644405365
$ws.Range("G7").Value = 1.544213533401489
$ws.Range("H7").Value = 0.3541924357414246
$ws.Range("I7").Value = 0.9842033386230469
$ws.Range("A8").Value = "model_6_2_4"
$ws.Range("B8").Value = 0.2933330629870743
$ws.Range("C8").Value = 0.1359172333389763
$ws.Range("D8").Value = -1.449571619219163
$ws.Range("E8").Value = 0.2528636580137626
$ws.Range("F8").Value = 0.7820720672607422
$ws.Range("G8").Value = 1.525166034698486
$ws.Range("H8").Value = 0.3585216104984283
$ws.Range("I8").Value = 0.9761565923690796
$ws.Range("A9").Value = "model_6_2_14"
$ws.Range("B9").Value = 0.29348756127548
$ws.Range("C9").Value = 0.1278102500403911
$ws.Range("D9").Value = -1.400731495095375
$ws.Range("E9").Value = 0.2496399596393868
$ws.Range("F9").Value = 0.7819010615348816
$ws.Range("G9").Value = 1.539475321769714
$ws.Range("H9").Value = 0.3513733148574829
$ws.Range("I9").Value = 0.980368435382843
$ws.Range("A10").Value = "model_6_2_10"
$ws.Range("B10").Value = 0.2943911627537656
$ws.Range("C10").Value = 0.1316191013542732
$ws.Range("D10").Value = -1.394662423753047
$ws.Range("E10").Value = 0.2526840655632612
$ws.Range("F10").Value = 0.7809010744094849
$ws.Range("G10").Value = 1.532752513885498
$ws.Range("H10").Value = 0.3504850566387177
$ws.Range("I10").Value = 0.9763913154602051
$ws.Range("A11").Value = "model_6_2_22"
$ws.Range("B11").Value = 0.2948009192297004
$ws.Range("C11").Value = 0.1326190484621492
$ws.Range("D11").Value = -1.40883226642906
$ws.Range("E11").Value = 0.2526523537437936
$ws.Range("F11").Value = 0.7804475426673889
$ws.Range("G11").Value = 1.53098738193512
$ws.Range("H11").Value = 0.3525589406490326
$ws.Range("I11").Value = 0.9764328002929688
$ws.Range("A12").Value = "model_6_2_5"
$ws.Range("B12").Value = 0.2949132700251683
$ws.Range("C12").Value = 0.136784885547595
$ws.Range("D12").Value = -1.404604702488317
$ws.Range("E12").Value = 0.2558544576519743
$ws.Range("F12").Value = 0.7803232073783875
$ws.Range("G12").Value = 1.523634433746338
$ws.Range("H12").Value = 0.3519402146339417
$ws.Range("I12").Value = 0.9722490906715393
$ws.Range("A13").Value = "model_6_2_16"
$ws.Range("B13").Value = 0.2949577797330135
$ws.Range("C13").Value = 0.1333532530024789
$ws.Range("D13").Value = -1.400148378658699
$ws.Range("E13").Value = 0.2536350278333619
$ws.Range("F13").Value = 0.780273973941803
$ws.Range("G13").Value = 1.529691457748413
$ws.Range("H13").Value = 0.3512879610061646
$ws.Range("I13").Value = 0.9751487970352173
$ws.Range("A14").Value = "model_6_2_17"
$ws.Range("B14").Value = 0.2949632537442488
$ws.Range("C14").Value = 0.133170275098434
$ws.Range("D14").Value = -1.397417648757259
$ws.Range("E14").Value = 0.2536480735664385
$ws.Range("F14").Value = 0.7802680134773254
$ws.Range("G14").Value = 1.530014514923096
$ws.Range("H14").Value = 0.3508882820606232
$ws.Range("I14").Value = 0.9751317501068115
$ws.Range("A15").Value = "model_6_2_18"
$ws.Range("B15").Value = 0.295121772602837
$ws.Range("C15").Value = 0.1330242846072452
$ws.Range("D15").Value = -1.389391699001949
$ws.Range("E15").Value = 0.2539669578122902
$ws.Range("F15").Value = 0.7800924181938171
$ws.Range("G15").Value = 1.530272126197815
$ws.Range("H15").Value = 0.3497136235237122
$ws.Range("I15").Value = 0.9747151732444763
$ws.Range("A16").Value = "model_6_2_15"
$ws.Range("B16").Value = 0.2957177567449119
$ws.Range("C16").Value = 0.1361873263049456
$ws.Range("D16").Value = -1.397527668325239
$ws.Range("E16").Value = 0.2558003602487877
$ws.Range("F16").Value = 0.7794329524040222
$ws.Range("G16").Value = 1.524689197540283
$ws.Range("H16").Value = 0.3509044051170349
$ws.Range("I16").Value = 0.9723197817802429
$ws.Range("A17").Value = "model_6_2_9"
$ws.Range("B17").Value = 0.2961513408142268
$ws.Range("C17").Value = 0.138179355227074
$ws.Range("D17").Value = -1.386423240294278
$ws.Range("E17").Value = 0.2578107726053988
$ws.Range("F17").Value = 0.7789530158042908
$ws.Range("G17").Value = 1.521173000335693
$ws.Range("H17").Value = 0.3492791652679443
$ws.Range("I17").Value = 0.969693124294281
$ws.Range("A18").Value = "model_6_2_13"
$ws.Range("B18").Value = 0.2963881874694364
$ws.Range("C18").Value = 0.138696750414979
$ws.Range("D18").Value = -1.393057196736398
$ws.Range("E18").Value = 0.257830841904193
$ws.Range("F18").Value = 0.7786909341812134
$ws.Range("G18").Value = 1.520259857177734
$ws.Range("H18").Value = 0.3502501249313354
$ws.Range("I18").Value = 0.9696668982505798
$ws.Range("A19").Value = "model_6_2_8"
$ws.Range("B19").Value = 0.2964163810694116
$ws.Range("C19").Value = 0.1407545324911047
$ws.Range("D19").Value = -1.399899116569243
$ws.Range("E19").Value = 0.2589416484010726
$ws.Range("F19").Value = 0.7786597013473511
$ws.Range("G19").Value = 1.516627788543701
$ws.Range("H19").Value = 0.351251482963562
$ws.Range("I19").Value = 0.9682155847549438
$ws.Range("A20").Value = "model_6_2_6"
$ws.Range("B20").Value = 0.2964484177040639
$ws.Range("C20").Value = 0.1406199027729237
$ws.Range("D20").Value = -1.389879497227462
$ws.Range("E20").Value = 0.2593738086823931
$ws.Range("F20").Value = 0.7786242365837097
$ws.Range("G20").Value = 1.516865253448486
$ws.Range("H20").Value = 0.3497850298881531
$ws.Range("I20").Value = 0.9676508903503418
$ws.Range("A21").Value = "model_6_2_21"
$ws.Range("B21").Value = 0.2968839050171601
$ws.Range("C21").Value = 0.140554545972241
$ws.Range("D21").Value = -1.404054898945989
$ws.Range("E21").Value = 0.2585797589431456
$ws.Range("F21").Value = 0.7781423330307007
$ws.Range("G21").Value = 1.516980648040771
$ws.Range("H21").Value = 0.3518597483634949
$ws.Range("I21").Value = 0.968688428401947
$ws.Range("A22").Value = "model_6_2_23"
$ws.Range("B22").Value = 0.2970281357903348
$ws.Range("C22").Value = 0.1405054720375397
$ws.Range("D22").Value = -1.400268348666207
$ws.Range("E22").Value = 0.2587441911039828
$ws.Range("F22").Value = 0.7779827117919922
$ws.Range("G22").Value = 1.517067432403564
$ws.Range("H22").Value = 0.3513055145740509
$ws.Range("I22").Value = 0.9684735536575317
$ws.Range("A23").Value = "model_6_2_19"
$ws.Range("B23").Value = 0.2983046408178568
$ws.Range("C23").Value = 0.145636399893712
$ws.Range("D23").Value = -1.392513318120053
$ws.Range("E23").Value = 0.2628228248463362
$ws.Range("F23").Value = 0.7765699625015259
$ws.Range("G23").Value = 1.508010983467102
$ws.Range("H23").Value = 0.3501704931259155
$ws.Range("I23").Value = 0.9631447792053223
$ws.Range("A24").Value = "model_6_2_11"
$ws.Range("B24").Value = 0.2989771953644103
$ws.Range("C24").Value = 0.1485809616048107
$ws.Range("D24").Value = -1.384392497963128
$ws.Range("E24").Value = 0.2653570745642803
$ws.Range("F24").Value = 0.7758256793022156
$ws.Range("G24").Value = 1.502813577651978
$ws.Range("H24").Value = 0.3489819169044495
$ws.Range("I24").Value = 0.9598336219787598
$ws.Range("A25").Value = "model_6_2_7"
$ws.Range("B25").Value = 0.2994351943494277
$ws.Range("C25").Value = 0.1505186698826096
$ws.Range("D25").Value = -1.37064978367047
$ws.Range("E25").Value = 0.2674672068363617
$ws.Range("F25").Value = 0.7753188014030457
$ws.Range("G25").Value = 1.499393224716187
$ws.Range("H25").Value = 0.3469705283641815
$ws.Range("I25").Value = 0.9570766687393188
$ws.Range("A26").Value = "model_6_2_24"
$ws.Range("B26").Value = 0.300025823579237
$ws.Range("C26").Value = 0.1526031866558233
$ws.Range("D26").Value = -1.408985738635839
$ws.Range("E26").Value = 0.266937029696071
$ws.Range("F26").Value = 0.7746652364730835
$ws.Range("G26").Value = 1.495713829994202
$ws.Range("H26").Value = 0.3525814116001129
$ws.Range("I26").Value = 0.9577692747116089